$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: fill in the sale-line data for the new transaction.
$ws.Range("A7").Value = 1

$ws.Range("C7").Value = "INSULINAGYPT 70/30 100 I.U./ML (4ML) VIAL"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "122.00"
$ws.Range("P7").Value = "122.0000"
$ws.Range("Q7").Value = "1:0"

$ws.Range("P8").Value = 122

# Footer timestamp refresh.
$ws.Range("A9").Value = "Tuesday, 16 September, 2025 9:42 AM"
